# Update the public EPEX Spot prices workbook:
# insert a new row for 2025-06-24 just before the existing 2025-06-25 row,
# on both the "Gaz" and "CO2" sheets (the "Prix Spot" sheet is untouched).

$wb = $excel.ActiveWorkbook

function Insert-PriceRow {
    param(
        [string]$SheetName,
        [int]$RowIndex,
        [string]$DateText,
        [double]$Price
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Push the existing row (and everything below it) down by one.
    $ws.Rows.Item($RowIndex).Insert()

    # Column A holds dates stored as plain text (e.g. "2025-06-25"), not
    # real Excel date serials, so force text formatting before assigning
    # the value - otherwise Excel's smart-entry would reinterpret the
    # string as a date. Reset the style back to Normal afterwards so the
    # new cell matches the unstyled look of its neighbours.
    $dateCell = $ws.Range("A" + $RowIndex)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $DateText
    $dateCell.Style = "Normal"

    $ws.Range("B" + $RowIndex).Value = $Price
}

# "Gaz" sheet: new row 10 = 2025-06-24 / 40.9, existing row shifts to 11.
Insert-PriceRow "Gaz" 10 "2025-06-24" 40.9

# "CO2" sheet: new row 10 = 2025-06-24 / 71.88, existing row shifts to 11.
Insert-PriceRow "CO2" 10 "2025-06-24" 71.88
